$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 4705.1333
$ws.Range("J32").Value = 5498.3335
$ws.Range("L32").Value = 5498.3335
$ws.Range("N32").Value = -6150.3335
$ws.Range("H62").Value = 5400.5713
$ws.Range("I62").Value = 5400.5713
$ws.Range("K62").Value = 5400.5713
$ws.Range("M62").Value = -4776.5713
$ws.Range("H64").Value = 6199.4
$ws.Range("J64").Value = 6199.4
$ws.Range("L64").Value = 6199.4
$ws.Range("N64").Value = -6695.4
$ws.Range("H65").Value = 5400.5713
$ws.Range("I65").Value = 5400.5713
$ws.Range("K65").Value = 27002.8565
$ws.Range("M65").Value = -23882.8565
$ws.Range("H67").Value = 6199.4
$ws.Range("J67").Value = 6199.4
$ws.Range("L67").Value = 6199.4
$ws.Range("N67").Value = -7915.4
$ws.Range("H74").Value = 5964
$ws.Range("I74").Value = 5874
$ws.Range("K74").Value = 5874
$ws.Range("M74").Value = -4938
$ws.Range("H77").Value = 5964
$ws.Range("I77").Value = 5874
$ws.Range("K77").Value = 29370
$ws.Range("M77").Value = -24690
$ws.Range("H98").Value = 15990.75
$ws.Range("I98").Value = 16573.018
$ws.Range("K98").Value = 16573.018
$ws.Range("M98").Value = -15075.018
$ws.Range("H107").Value = 783.9
$ws.Range("J107").Value = 567.8
$ws.Range("L107").Value = 567.8
$ws.Range("N107").Value = -4407.8
$ws.Range("H116").Value = 8517.632
$ws.Range("I116").Value = 7340.5
$ws.Range("J116").Value = 10535.571
$ws.Range("K116").Value = 7340.5
$ws.Range("L116").Value = 10535.571
$ws.Range("M116").Value = -3898.5
$ws.Range("N116").Value = -17419.571
$ws.Range("H122").Value = 15990.75
$ws.Range("I122").Value = 16573.018
$ws.Range("K122").Value = 49719.054
$ws.Range("M122").Value = -47269.054
$ws.Range("H137").Value = 54058916
$ws.Range("I137").Value = 43481996
$ws.Range("J137").Value = 71435290
$ws.Range("K137").Value = 130445988
$ws.Range("L137").Value = 214305870
$ws.Range("M137").Value = -130443438
$ws.Range("N137").Value = -214310970
$ws.Range("H138").Value = 4228492.5
$ws.Range("J138").Value = 4704386
$ws.Range("L138").Value = 14113158
$ws.Range("N138").Value = -14123438

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18193022
$ws.Range("I32").Value = 20416618
$ws.Range("K32").Value = 20416618
$ws.Range("M32").Value = -20416331
$ws.Range("H45").Value = 2288.9167
$ws.Range("J45").Value = 2413.6
$ws.Range("L45").Value = 2413.6
$ws.Range("N45").Value = -3167.6
$ws.Range("H61").Value = 19611904
$ws.Range("I61").Value = 22731026
$ws.Range("J61").Value = 5999.2856
$ws.Range("K61").Value = 22731026
$ws.Range("L61").Value = 5999.2856
$ws.Range("M61").Value = -22730814
$ws.Range("N61").Value = -6423.2856
$ws.Range("H74").Value = 60607616
$ws.Range("I74").Value = 74075384
$ws.Range("K74").Value = 74075384
$ws.Range("M74").Value = -74074510
$ws.Range("H77").Value = 60607616
$ws.Range("I77").Value = 74075384
$ws.Range("K77").Value = 370376920
$ws.Range("M77").Value = -370372552
$ws.Range("H132").Value = 31259346
$ws.Range("I132").Value = 10631
$ws.Range("K132").Value = 31893
$ws.Range("M132").Value = -29363
$ws.Range("H136").Value = 19611904
$ws.Range("I136").Value = 22731026
$ws.Range("J136").Value = 5999.2856
$ws.Range("K136").Value = 68193078
$ws.Range("L136").Value = 17997.8568
$ws.Range("M136").Value = -68190528
$ws.Range("N136").Value = -23097.8568

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3798.8333
$ws.Range("I94").Value = 4156.6
$ws.Range("J94").Value = 2010
$ws.Range("K94").Value = 4156.6
$ws.Range("L94").Value = 2010
$ws.Range("M94").Value = -3705.6
$ws.Range("N94").Value = -2912
$ws.Range("H107").Value = 6385.5835
$ws.Range("I107").Value = 5864.1113
$ws.Range("K107").Value = 5864.1113
$ws.Range("M107").Value = -3944.1113
$ws.Range("H134").Value = 2107.3914
$ws.Range("I134").Value = 2404.8667
$ws.Range("J134").Value = 1549.625
$ws.Range("K134").Value = 7214.6001
$ws.Range("L134").Value = 4648.875
$ws.Range("M134").Value = -4679.6001
$ws.Range("N134").Value = -9718.875

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 20413940
$ws.Range("I31").Value = 4103.769
$ws.Range("J31").Value = 43485930
$ws.Range("K31").Value = 4103.769
$ws.Range("L31").Value = 43485930
$ws.Range("M31").Value = -3808.769
$ws.Range("N31").Value = -43486520
$ws.Range("H34").Value = 20413940
$ws.Range("I34").Value = 4103.769
$ws.Range("J34").Value = 43485930
$ws.Range("K34").Value = 4103.769
$ws.Range("L34").Value = 43485930
$ws.Range("M34").Value = -3901.769
$ws.Range("N34").Value = -43486334
$ws.Range("H62").Value = 5549.6665
$ws.Range("I62").Value = 3499
$ws.Range("K62").Value = 3499
$ws.Range("M62").Value = -2875
$ws.Range("H65").Value = 5549.6665
$ws.Range("I65").Value = 3499
$ws.Range("K65").Value = 17495
$ws.Range("M65").Value = -14375

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 175
$ws.Range("I17").Value = 150
$ws.Range("K17").Value = 450
$ws.Range("M17").Value = -281
$ws.Range("H76").Value = 3547.8
$ws.Range("I76").Value = 2681.3333
$ws.Range("J76").Value = 4847.5
$ws.Range("K76").Value = 8043.999899999999
$ws.Range("L76").Value = 14542.5
$ws.Range("M76").Value = -7660.999899999999
$ws.Range("N76").Value = -15308.5
$ws.Range("H79").Value = 3547.8
$ws.Range("I79").Value = 2681.3333
$ws.Range("J79").Value = 4847.5
$ws.Range("K79").Value = 8043.999899999999
$ws.Range("L79").Value = 14542.5
$ws.Range("M79").Value = -6717.999899999999
$ws.Range("N79").Value = -17194.5
$ws.Range("H82").Value = 4432
$ws.Range("I82").Value = 3318.6
$ws.Range("J82").Value = 9999
$ws.Range("K82").Value = 9955.799999999999
$ws.Range("L82").Value = 29997
$ws.Range("M82").Value = -9549.799999999999
$ws.Range("N82").Value = -30809
$ws.Range("H85").Value = 4432
$ws.Range("I85").Value = 3318.6
$ws.Range("J85").Value = 9999
$ws.Range("K85").Value = 9955.799999999999
$ws.Range("L85").Value = 29997
$ws.Range("M85").Value = -8551.799999999999
$ws.Range("N85").Value = -32805
$ws.Range("H131").Value = 1865.0769
$ws.Range("J131").Value = 1846.6
$ws.Range("L131").Value = 5539.799999999999
$ws.Range("N131").Value = -15619.8
$ws.Range("H132").Value = 3178762.2
$ws.Range("J132").Value = 3708111.5
$ws.Range("L132").Value = 33373003.5
$ws.Range("N132").Value = -33378063.5

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H108").Value = 85000
$ws.Range("J108").Value = 85000
$ws.Range("L108").Value = 85000
$ws.Range("N108").Value = -92680
$ws.Range("H132").Value = 5628.763
$ws.Range("I132").Value = 4950.241
$ws.Range("J132").Value = 7815.1113
$ws.Range("K132").Value = 14850.723
$ws.Range("L132").Value = 23445.3339
$ws.Range("M132").Value = -12320.723
$ws.Range("N132").Value = -28505.3339

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1224.7
$ws.Range("I16").Value = 1163.4445
$ws.Range("J16").Value = 1776
$ws.Range("K16").Value = 1163.4445
$ws.Range("L16").Value = 1776
$ws.Range("M16").Value = -993.4445000000001
$ws.Range("N16").Value = -2116
$ws.Range("H46").Value = 1716.3
$ws.Range("I46").Value = 795.06665
$ws.Range("K46").Value = 795.06665
$ws.Range("M46").Value = -607.06665
$ws.Range("H68").Value = 4754.9375
$ws.Range("I68").Value = 3331.6667
$ws.Range("J68").Value = 5083.385
$ws.Range("K68").Value = 3331.6667
$ws.Range("L68").Value = 5083.385
$ws.Range("M68").Value = -2582.6667
$ws.Range("N68").Value = -6581.385
$ws.Range("H71").Value = 4754.9375
$ws.Range("I71").Value = 3331.6667
$ws.Range("J71").Value = 5083.385
$ws.Range("K71").Value = 16658.3335
$ws.Range("L71").Value = 25416.925
$ws.Range("M71").Value = -12914.3335
$ws.Range("N71").Value = -32904.925
$ws.Range("H100").Value = 4317.8335
$ws.Range("I100").Value = 3604
$ws.Range("J100").Value = 4460.6
$ws.Range("K100").Value = 3604
$ws.Range("L100").Value = 4460.6
$ws.Range("M100").Value = -3063
$ws.Range("N100").Value = -5542.6
$ws.Range("H135").Value = 87040.86
$ws.Range("J135").Value = 87040.86
$ws.Range("L135").Value = 87040.86
$ws.Range("N135").Value = -97180.86

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()
$ws.Range("H113").Value = 782.5714
$ws.Range("I113").Value = 578.7778
$ws.Range("J113").Value = 1149.4
$ws.Range("K113").Value = 1736.3334
$ws.Range("L113").Value = 3448.2
$ws.Range("M113").Value = 433.6666
$ws.Range("N113").Value = -7788.200000000001
$ws.Range("H136").Value = 1317.2046
$ws.Range("I136").Value = 1326.3
$ws.Range("K136").Value = 3978.9
$ws.Range("M136").Value = -1428.9
